# "Added '[Wiki]' to the front of Wiki page names in the results display.
#  Added year to the front of forum subjects in the results display."
#
# The author illustrated this J Viewer app-behaviour change by dropping a
# screenshot of the updated results list onto the info slide (slide 2 /
# sldId 257), alongside the existing White Rabbit artwork. That shows up in
# the collab-tracking metadata as a new "add" picChg (picMk id="2") plus a
# refreshed "mod" picChg on the pre-existing picture (the rabbit, id 1026)
# sitting next to it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$imagePath = "/tmp/work/results_screenshot.png"

# Drop it in under the rabbit artwork (points, like real PPT COM), without
# touching any of the existing shapes' size/position.
$left   = 700
$top    = 290
$width  = 230
$height = 153

$pic = $s.Shapes.AddPicture($imagePath, $false, $true, $left, $top, $width, $height)
$pic.Name = "Picture 3"
